# SignIn function has been added
# Updates the RegistrationData / LoginData test-data workbook:
#  - refresh a stale test email address (eddie309 -> eddie510)
#  - fix a typo'd login test row (gamil.com -> gmail.com) and refresh its password
#  - drop two obsolete login test rows
#  - make LoginData the active sheet/selection, matching the saved UI state

$wb  = $excel.ActiveWorkbook
$wsReg   = $wb.Worksheets.Item("RegistrationData")
$wsLogin = $wb.Worksheets.Item("LoginData")

# --- RegistrationData: refresh the test email address used for Eddie/He ---
$wsReg.Cells.Item(2, 3).Value = "eddie510@gmail.com"

# --- LoginData: update row 2 (email + password) and fix row 3's email typo ---
$wsLogin.Cells.Item(2, 1).Value = "eddie510@gmail.com"
$wsLogin.Cells.Item(2, 2).Value = "eddie123"
$wsLogin.Cells.Item(3, 1).Value = "eddie123@gmail.com"

# --- LoginData: remove the two trailing obsolete rows (old rows 9 and 10) ---
$wsLogin.Rows.Item(9).Delete()
$wsLogin.Rows.Item(9).Delete()

# --- Restore selections and make LoginData the active/selected sheet ---
$wsReg.Range("I14").Select()

$wsLogin.Activate()
$wsLogin.Range("A20").Select()
